$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8
$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44911
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100112039
$ws.Cells.Item($row, 7).Value = "Ciboulette"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 700
$ws.Cells.Item($row, 11).Value = 1800
$ws.Cells.Item($row, 12).Value = 2000
$ws.Cells.Item($row, 13).Value = 1900
$ws.Cells.Item($row, 14).Value = "$/docena de atados"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 633
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
